# Update NATMI LR-pair edge-weight sheet with newly recomputed TPM-based
# ligand/receptor expression statistics ("update scripts wuth new tpm").
#
# The sheet has one row per (sending cluster, target cluster) combination.
# Columns G/H/I/J (ligand avg / total expression + derived specificities)
# depend only on the sending cluster (column A); columns M/N/O/P (receptor
# avg / total expression + derived specificities) depend only on the
# target cluster (column D). Columns Q/R (edge avg/total expression
# weight) are simply G*M and H*N, and S/T are those edge weights
# normalised (derived specificity) across every row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1
$firstDataRow = 2

# New, recomputed ligand-side values (TPM-based) keyed by sending cluster name.
$ligandStats = @{
    "ECs"               = @(2.0407475,            4.081495,    0.007689531190315315, 0.005145434208836267)
    "FAPs"              = @(108.0898166666667,     324.26945,   0.4072821437310581,   0.4087980313366845)
    "Inflammatory-Mac"  = @(49.59263000000001,     148.77789,   0.18686489886415,     0.1875604024320694)
    "MuSCs"             = @(0.9116095,             1.823219,    0.003434942188407801, 0.002298484602529281)
    "Neutrophils"       = @(97.62255466666666,     292.867664,  0.3678415281594588,   0.3692106193949926)
    "Resolving-Mac"     = @(7.135609000000001,     21.406827,   0.02688695586661,     0.0269870280248879)
}

# New, recomputed receptor-side values (TPM-based) keyed by target cluster name.
$receptorStats = @{
    "Inflammatory-Mac" = @(7.318981333333333, 21.956944,           0.1531761772116645, 0.1531761772116645)
    "Neutrophils"      = @(32.599203,         97.79760900000001,   0.6822563234237459, 0.6822563234237459)
    "Resolving-Mac"    = @(7.863275333333333, 23.589826,           0.1645674993645896, 0.1645674993645896)
}

# Pass 1: write ligand (G,H,I,J) and receptor (M,N,O,P) columns, and
# compute the raw (un-normalised) edge weights Q = G*M and R = H*N.
$edgeAvg = @{}
$edgeTotal = @{}
$sumEdgeAvg = 0
$sumEdgeTotal = 0

for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $sender = $ws.Cells.Item($r, 1).Value2
    $target = $ws.Cells.Item($r, 4).Value2
    if (-not $sender -or -not $target) { continue }

    $lig = $ligandStats[$sender]
    $rec = $receptorStats[$target]
    if (-not $lig -or -not $rec) { continue }

    $ws.Cells.Item($r, 7).Value  = $lig[0]   # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $lig[1]   # H: Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $lig[2]   # I: Ligand derived specificity (avg)
    $ws.Cells.Item($r, 10).Value = $lig[3]   # J: Ligand derived specificity (total)

    $ws.Cells.Item($r, 13).Value = $rec[0]   # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $rec[1]   # N: Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $rec[2]   # O: Receptor derived specificity (avg)
    $ws.Cells.Item($r, 16).Value = $rec[3]   # P: Receptor derived specificity (total)

    $q = $lig[0] * $rec[0]
    $rr = $lig[1] * $rec[1]
    $edgeAvg[$r] = $q
    $edgeTotal[$r] = $rr
    $sumEdgeAvg = $sumEdgeAvg + $q
    $sumEdgeTotal = $sumEdgeTotal + $rr
}

# Pass 2: write Q,R (edge weights) and S,T (edge weights normalised over
# the whole sheet -> "derived specificity" of the edge weight).
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    if (-not $edgeAvg.ContainsKey($r)) { continue }

    $q = $edgeAvg[$r]
    $rr = $edgeTotal[$r]

    $ws.Cells.Item($r, 17).Value = $q                          # Q
    $ws.Cells.Item($r, 18).Value = $rr                         # R
    $ws.Cells.Item($r, 19).Value = $q / $sumEdgeAvg            # S
    $ws.Cells.Item($r, 20).Value = $rr / $sumEdgeTotal         # T
}

Write-Output "Updated rows $firstDataRow..$lastRow with new TPM-derived statistics."
